# Weekly price-report update: a new observation (price record) for
# "Pepino ensalada" at Vega Monumental Concepción was inserted as the
# first data row (row 13), pushing all subsequent rows down by one.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 13, shifting rows 13:103 down
# to 14:104 (this also grows the used range to A1:R104, matching the
# other rows' D-column date style for the new row).
$ws.Rows.Item(13).Insert()

# Populate the newly inserted row 13 with the new record's data. The
# descriptive columns (market id/name/region/category id/category/
# variety/quality/classification) are constant for every row in this
# sheet, so reuse the same values as the surrounding rows.
$ws.Range("A13").Value2 = 11
$ws.Range("B13").Value2 = "Vega Monumental Concepción"
$ws.Range("C13").Value2 = "Bíobío"
$ws.Range("D13").Value2 = 44602
$ws.Range("E13").Value2 = 8
$ws.Range("F13").Value2 = 100112043
$ws.Range("G13").Value2 = "Pepino ensalada"
$ws.Range("H13").Value2 = "Sin especificar"
$ws.Range("I13").Value2 = "Primera"
$ws.Range("J13").Value2 = 220
$ws.Range("K13").Value2 = 10000
$ws.Range("L13").Value2 = 11000
$ws.Range("M13").Value2 = 10545
$ws.Range("N13").Value2 = "`$/caja 80 unidades"
$ws.Range("O13").Value2 = "Región del Maule"
$ws.Range("P13").Value2 = 132
$ws.Range("Q13").Value2 = 80
$ws.Range("R13").Value2 = "Hortaliza"
